$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the trailing space from the "Baz " shared string used by A2.
$ws.Cells.Item(2, 1).Value2 = "Baz"

# Give B2 a thin box border on all four sides (adds a new border/cellXf).
$ws.Cells.Item(2, 2).Borders.LineStyle = 1

# Row 2's height was nudged from 15 to 14.9.
$ws.Rows.Item(2).RowHeight = 14.9

# Selection moved from B1 to B2.
$ws.Range("B2").Select() | Out-Null
